$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "89.177.55"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +9.77%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.385.68"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +7.63%  "
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "223.41"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.97%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "654.43"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.33%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.426"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +50.45%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.670"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +15.02%  "
$ws.Range("E9").Value = "  -0.02%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.379.03"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("E11").Value = "  +9.30%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000293"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +16.41%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "37.32"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +18.79%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.170"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.08%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.67"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +7.40%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "4.008.30"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +7.68%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "89.012.56"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +9.75%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.373.41"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +7.48%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "15.03"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +8.14%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "469.72"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +8.65%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.75"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +9.39%  "
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.95%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.74"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +12.50%  "
$ws.Range("E24").Value = "  +3.98%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.62"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.46%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.91"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +19.53%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "3.561.24"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +7.63%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0000146"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +20.86%  "
$ws.Range("B29").Value = "Litecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "81.52"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.98%  "
$ws.Range("E30").Value = "  +0.03%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.195"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +40.54%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "9.52"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +6.31%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "595.40"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +5.19%  "
$ws.Range("E36").Value = "  +7.41%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.38"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +22.90%  "
$ws.Range("E38").Value = "  -5.62%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "24.07"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +6.01%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.440"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +8.14%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +6.22%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.25"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.39%  "
$ws.Range("E43").Value = "  +5.47%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "158.68"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +10.48%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "190.22"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "46.96"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.42%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.58"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +9.49%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.678"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +8.62%  "
